$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = $ws.Range("F3:G9")
$target.NumberFormat = "0.00"
$target.HorizontalAlignment = -4108

$ws.Range("F3").Formula = "=C3*3.28084"
$ws.Range("G3").Formula = "=D3*3.28084"
$ws.Range("F4:F9").Formula = "=C4*3.28084"
$ws.Range("G4:G9").Formula = "=D4*3.28084"

Write-Host "done"
